$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add dates to B13, B14, B15, B16 matching the style used elsewhere in column B (style index 1 => numFmt 16)
$ws.Range("B13").Value = 41721
$ws.Range("B14").Value = 41721
$ws.Range("B15").Value = 41723
$ws.Range("B16").Value = 41724

$ws.Range("B13:B16").NumberFormat = $ws.Range("B12").NumberFormat

# Update active selection to H19
$ws.Range("H19").Select()
